# Auto-generated Excel COM-interop script to apply F-column numeric updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1666
$ws.Range("F7").Value = 314
$ws.Range("F9").Value = 3505
$ws.Range("F10").Value = 919
$ws.Range("F11").Value = 1152
$ws.Range("F16").Value = 1254
$ws.Range("F17").Value = 1781
$ws.Range("F20").Value = 1540
$ws.Range("F21").Value = 1071
$ws.Range("F22").Value = 2085
$ws.Range("F23").Value = 143
$ws.Range("F24").Value = 4243
$ws.Range("F26").Value = 2693
$ws.Range("F27").Value = 1210

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 163
$ws.Range("F14").Value = 34
$ws.Range("F15").Value = 34
$ws.Range("F20").Value = 13
$ws.Range("F23").Value = 115
$ws.Range("F41").Value = 15
$ws.Range("F44").Value = 82

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2548
$ws.Range("F6").Value = 9571
$ws.Range("F11").Value = 2953
$ws.Range("F12").Value = 460
$ws.Range("F13").Value = 791
$ws.Range("F14").Value = 192

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 2953
$ws.Range("F9").Value = 791
$ws.Range("F13").Value = 314
$ws.Range("F16").Value = 919
$ws.Range("F17").Value = 1152
$ws.Range("F23").Value = 1254
$ws.Range("F26").Value = 34
$ws.Range("F27").Value = 34
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 1781
$ws.Range("F33").Value = 1540
$ws.Range("F35").Value = 115
$ws.Range("F36").Value = 115
$ws.Range("F38").Value = 1071
$ws.Range("F40").Value = 2085
$ws.Range("F42").Value = 143
$ws.Range("F44").Value = 4243
$ws.Range("F46").Value = 2693
